$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) The cached "datetimeFigureOut" field text on the slide master and on all
#    11 slide layouts rolled from 11/11/2018 to 12/11/2018 (PowerPoint
#    refreshes these automatic date fields whenever the deck is saved on a
#    later day).
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "11/11/2018") {
                $sh.TextFrame.TextRange.Text = "12/11/2018"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------------------
# 2) Reposition the "add ... Modules(m, index)" label text box slightly
#    (refactor of the events containing "Module").
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 217") {
        $sh.Left = 225.917724609375
        $sh.Top = 366.9625244140625
    }
}

# ---------------------------------------------------------------------------
# 3) Add a new dashed connector line mirroring "Straight Connector 70",
#    covering the extra "Module" lifeline span that was introduced.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Straight Connector 70") {
        $newConn = $sh.Duplicate().Item(1)
        $newConn.Name = "Straight Connector 71"
        $newConn.Left = 247.69677734375
        $newConn.Top = 293.8426208496094
        $newConn.Width = 0.8771653771400452
        $newConn.Height = 136.492919921875
    }
}
